# Apply the edits described by the diff to sheet "Hoja1" of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# A1: date serial 45406 (2024-04-24) -> 45436 (2024-05-24)
$ws.Range("A1").Value = 45436

# D29: 8541 -> 23711
$ws.Range("D29").Value = 23711

# D30: 10272 -> 26200
$ws.Range("D30").Value = 26200
